# Stock on hand - rebuild the "Data" table with the new column layout:
#  - remove Cost_Amount
#  - remove the Pre/In/Post quantity columns
#  - add Stock_On_Hand (moved up), Total_Cost, Quantity_Purchases, Total_Purchases,
#    Quantity_Consumed, Total_Consumption
#  - move Location_Code to the end (just before StartDate/EndDate)
#  - restyle the header (light-grey fill, left aligned, text format) and the
#    data row (centered, not bold)
#  - switch the table look to TableStyleMedium6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Tear down the existing table but keep the sheet around ------------
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()
$ws.Range("A1:Q2").Clear()

# --- 2. New header / data layout ------------------------------------------
$headers = @(
    "Item_No",
    "Description",
    "Unit_Cost",
    "Stock_On_Hand",
    "Total_Cost",
    "Quantity_Purchases",
    "Total_Purchases",
    "Quantity_Consumed",
    "Total_Consumption",
    "BaseUnit",
    "Gen_Prod_Posting_Group",
    "Item_Catogry",
    "Location_Code",
    "StartDate",
    "EndDate"
)

# column kinds: "text" -> blank string data row, "num" -> 0, "date" -> 0 (date fmt)
$kinds = @(
    "text", "text", "num", "num", "num", "num", "num", "num", "num",
    "text", "text", "text", "text", "date", "date"
)

$colCount = $headers.Count

for ($c = 1; $c -le $colCount; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
    $kind = $kinds[$c - 1]
    if ($kind -eq "num") {
        $ws.Cells.Item(2, $c).Value = 0
    } elseif ($kind -eq "date") {
        $ws.Cells.Item(2, $c).Value = 0
    } else {
        $ws.Cells.Item(2, $c).Value = ""
    }
}

$lastColLetter = "O"
$headerRange = $ws.Range("A1:" + $lastColLetter + "1")
$dataRange = $ws.Range("A2:" + $lastColLetter + "2")
$fullRange = $ws.Range("A1:" + $lastColLetter + "2")

# --- 3. Recreate the table --------------------------------------------------
$newLo = $ws.ListObjects.Add(1, $fullRange, $null, 1)
$newLo.Name = "Data"
$newLo.TableStyle = "TableStyleMedium6"

# --- 4. Header formatting ---------------------------------------------------
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Amasis MT Pro Light"
$headerRange.Font.Size = 12
$headerRange.Interior.Color = 192 + 192 * 256 + 192 * 65536
$headerRange.HorizontalAlignment = -4131   # xlLeft
$headerRange.VerticalAlignment = -4108     # xlCenter
$headerRange.NumberFormat = "@"
$ws.Rows.Item(1).RowHeight = 15.6

# --- 5. Data row formatting --------------------------------------------------
$dataRange.Font.Bold = $false
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 11
$dataRange.HorizontalAlignment = -4108     # xlCenter
$dataRange.VerticalAlignment = -4108       # xlCenter

foreach ($addr in @("A2", "B2", "J2", "K2", "L2", "M2")) {
    $ws.Range($addr).NumberFormat = "General"
}
foreach ($addr in @("C2", "D2", "E2", "F2", "G2", "H2", "I2")) {
    $ws.Range($addr).NumberFormat = "#,##0.00"
}
foreach ($addr in @("N2", "O2")) {
    $ws.Range($addr).NumberFormat = "m/d/yyyy"
}

# --- 6. Misc sheet view bits -------------------------------------------------
$ws.Range("D4").Select()
